$wb = $excel.ActiveWorkbook

# Sheet ALC
$ws = $wb.Worksheets.Item(1)
# Row 13
$ws.Range("H13").Value = 50000
$ws.Range("I13").Value = 50000
$ws.Range("J13").Value = 50000
$ws.Range("K13").Value = 50000
$ws.Range("L13").Value = 50000
$ws.Range("M13").Value = -49831
$ws.Range("N13").Value = -50338

# Row 41
$ws.Range("H41").Value = 5848273
$ws.Range("I41").Value = 6944749
$ws.Range("J41").Value = 400
$ws.Range("K41").Value = 6944749
$ws.Range("L41").Value = 400
$ws.Range("M41").Value = -6944309
$ws.Range("N41").Value = -1280

# Row 98
$ws.Range("H98").Value = 431645
$ws.Range("I98").Value = 448750.8
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 448750.8
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = -447252.8
$ws.Range("N98").Value = -6996

# Row 121
$ws.Range("H121").Value = 618.8049
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 618.8049
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 1856.4147
$ws.Range("N121").Value = -5350.414699999999
$ws.Range("M121").ClearContents()

# Row 122
$ws.Range("H122").Value = 431645
$ws.Range("I122").Value = 448750.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 1346252.4
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1343802.4
$ws.Range("N122").Value = -16900

# Row 132
$ws.Range("H132").Value = 26598.44
$ws.Range("I132").Value = 27936.82
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 83810.45999999999
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -81280.45999999999
$ws.Range("N132").Value = -6560

# Row 133
$ws.Range("H133").Value = 49477.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 49477.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 49477.5
$ws.Range("N133").Value = -59597.5

# Row 138
$ws.Range("H138").Value = 6283063.5
$ws.Range("I138").Value = 3249538.8
$ws.Range("J138").Value = 7250274.5
$ws.Range("K138").Value = 9748616.399999999
$ws.Range("L138").Value = 21750823.5
$ws.Range("M138").Value = -9743476.399999999
$ws.Range("N138").Value = -21761103.5

# Row 141
$ws.Range("H141").Value = 3194.5
$ws.Range("I141").Value = 3216.111
$ws.Range("J141").Value = 3000
$ws.Range("K141").Value = 9648.332999999999
$ws.Range("L141").Value = 9000
$ws.Range("M141").Value = -4468.332999999999
$ws.Range("N141").Value = -19360

# Sheet ARM
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 36097.387
$ws.Range("I32").Value = 9388.289000000001
$ws.Range("J32").Value = 128365.18
$ws.Range("K32").Value = 9388.289000000001
$ws.Range("L32").Value = 128365.18
$ws.Range("M32").Value = -9101.289000000001
$ws.Range("N32").Value = -128939.18

# Row 61
$ws.Range("H61").Value = 2749.8333
$ws.Range("I61").Value = 2140.8
$ws.Range("J61").Value = 3967.9
$ws.Range("K61").Value = 2140.8
$ws.Range("L61").Value = 3967.9
$ws.Range("M61").Value = -1928.8
$ws.Range("N61").Value = -4391.9

# Row 74
$ws.Range("H74").Value = 5236.1816
$ws.Range("I74").Value = 1049.091
$ws.Range("J74").Value = 13610.363
$ws.Range("K74").Value = 1049.091
$ws.Range("L74").Value = 13610.363
$ws.Range("M74").Value = -175.0909999999999
$ws.Range("N74").Value = -15358.363

# Row 77
$ws.Range("H77").Value = 5236.1816
$ws.Range("I77").Value = 1049.091
$ws.Range("J77").Value = 13610.363
$ws.Range("K77").Value = 5245.455
$ws.Range("L77").Value = 68051.815
$ws.Range("M77").Value = -877.4549999999999
$ws.Range("N77").Value = -76787.815

# Row 133
$ws.Range("H133").Value = 39602.332
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 39602.332
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 39602.332
$ws.Range("N133").Value = -44662.332

# Row 136
$ws.Range("H136").Value = 2749.8333
$ws.Range("I136").Value = 2140.8
$ws.Range("J136").Value = 3967.9
$ws.Range("K136").Value = 6422.400000000001
$ws.Range("L136").Value = 11903.7
$ws.Range("M136").Value = -3872.400000000001
$ws.Range("N136").Value = -17003.7

# Sheet BSM
$ws = $wb.Worksheets.Item(3)
# Row 105
$ws.Range("H105").Value = 2705.875
$ws.Range("I105").Value = 2575.6177
$ws.Range("J105").Value = 3022.2144
$ws.Range("K105").Value = 2575.6177
$ws.Range("L105").Value = 3022.2144
$ws.Range("M105").Value = -828.6176999999998
$ws.Range("N105").Value = -6516.2144

# Sheet CRP
$ws = $wb.Worksheets.Item(4)
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

# Row 31
$ws.Range("H31").Value = 6861.175
$ws.Range("I31").Value = 6800
$ws.Range("J31").Value = 6866.1353
$ws.Range("K31").Value = 6800
$ws.Range("L31").Value = 6866.1353
$ws.Range("M31").Value = -6505
$ws.Range("N31").Value = -7456.1353

# Row 34
$ws.Range("H34").Value = 6861.175
$ws.Range("I34").Value = 6800
$ws.Range("J34").Value = 6866.1353
$ws.Range("K34").Value = 6800
$ws.Range("L34").Value = 6866.1353
$ws.Range("M34").Value = -6598
$ws.Range("N34").Value = -7270.1353

# Row 69
$ws.Range("H69").Value = 20000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498

# Row 72
$ws.Range("H72").Value = 20000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488

# Sheet CUL
$ws = $wb.Worksheets.Item(5)
# Row 3
$ws.Range("H3").Value = 3806.8823
$ws.Range("I3").Value = 1091.7
$ws.Range("J3").Value = 7685.7144
$ws.Range("K3").Value = 3275.1
$ws.Range("L3").Value = 23057.1432
$ws.Range("M3").Value = -3163.1
$ws.Range("N3").Value = -23281.1432

# Row 34
$ws.Range("H34").Value = 2800
$ws.Range("I34").Value = 400
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1200
$ws.Range("L34").Value = 12000
$ws.Range("M34").Value = -1116
$ws.Range("N34").Value = -12168

# Row 39
$ws.Range("H39").Value = 9263.315000000001
$ws.Range("I39").Value = 999
$ws.Range("J39").Value = 9722.444
$ws.Range("K39").Value = 2997
$ws.Range("L39").Value = 29167.332
$ws.Range("M39").Value = -2703
$ws.Range("N39").Value = -29755.332

# Row 55
$ws.Range("H55").Value = 2600
$ws.Range("I55").Value = 666.6667
$ws.Range("J55").Value = 5500
$ws.Range("K55").Value = 2000.0001
$ws.Range("L55").Value = 16500
$ws.Range("M55").Value = -1823.0001
$ws.Range("N55").Value = -16854

# Row 134
$ws.Range("H134").Value = 6903.6284
$ws.Range("I134").Value = 3324.889
$ws.Range("J134").Value = 10692.883
$ws.Range("K134").Value = 9974.667000000001
$ws.Range("L134").Value = 32078.649
$ws.Range("M134").Value = -4904.667000000001
$ws.Range("N134").Value = -42218.649

# Row 139
$ws.Range("H139").Value = 2522.8462
$ws.Range("I139").Value = 2522.8462
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 7568.5386
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = -2428.5386

# Sheet GSM
$ws = $wb.Worksheets.Item(6)
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# Row 122
$ws.Range("H122").Value = 2661.5356
$ws.Range("I122").Value = 2688.4
$ws.Range("J122").Value = 2437.6667
$ws.Range("K122").Value = 8065.200000000001
$ws.Range("L122").Value = 7313.000100000001
$ws.Range("M122").Value = -5615.200000000001
$ws.Range("N122").Value = -12213.0001

# Row 126
$ws.Range("H126").Value = 2642.1914
$ws.Range("I126").Value = 2351
$ws.Range("J126").Value = 2839.7856
$ws.Range("K126").Value = 7053
$ws.Range("L126").Value = 8519.356800000001
$ws.Range("M126").Value = -4583
$ws.Range("N126").Value = -13459.3568

# Row 132
$ws.Range("H132").Value = 2420.9167
$ws.Range("I132").Value = 2621.7585
$ws.Range("J132").Value = 1588.8572
$ws.Range("K132").Value = 7865.2755
$ws.Range("L132").Value = 4766.571599999999
$ws.Range("M132").Value = -5335.2755
$ws.Range("N132").Value = -9826.571599999999

# Sheet WVR
$ws = $wb.Worksheets.Item(8)
# Row 136
$ws.Range("H136").Value = 2540.6
$ws.Range("I136").Value = 790.26086
$ws.Range("J136").Value = 5895.4165
$ws.Range("K136").Value = 2370.78258
$ws.Range("L136").Value = 17686.2495
$ws.Range("M136").Value = 179.2174199999999
$ws.Range("N136").Value = -22786.2495
